# Auto-generated Excel COM-interop script to apply market price/profit
# data updates across multiple worksheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16835.166
$ws.Range("J21").Value = 19998.75
$ws.Range("L21").Value = 19998.75
$ws.Range("N21").Value = -20934.75
$ws.Range("H23").Value = 16835.166
$ws.Range("J23").Value = 19998.75
$ws.Range("L23").Value = 19998.75
$ws.Range("N23").Value = -20466.75
$ws.Range("H76").Value = 3354.5454
$ws.Range("I76").Value = 3290
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3290
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -2975
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 3354.5454
$ws.Range("I79").Value = 3290
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3290
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2198
$ws.Range("N79").Value = -6184
$ws.Range("H106").Value = 35295856
$ws.Range("I106").Value = 40001596
$ws.Range("K106").Value = 40001596
$ws.Range("M106").Value = -40000965
$ws.Range("H132").Value = 4041.4324
$ws.Range("I132").Value = 3967
$ws.Range("J132").Value = 4311.25
$ws.Range("K132").Value = 11901
$ws.Range("L132").Value = 12933.75
$ws.Range("M132").Value = -9371
$ws.Range("N132").Value = -17993.75
$ws.Range("H138").Value = 2584.5632
$ws.Range("I138").Value = 3374.6667
$ws.Range("J138").Value = 2458.1467
$ws.Range("K138").Value = 10124.0001
$ws.Range("L138").Value = 7374.4401
$ws.Range("M138").Value = -4984.000100000001
$ws.Range("N138").Value = -17654.4401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2138.9375
$ws.Range("I2").Value = 2226.9167
$ws.Range("K2").Value = 2226.9167
$ws.Range("M2").Value = -2113.9167
$ws.Range("H32").Value = 372519.38
$ws.Range("I32").Value = 407234.9
$ws.Range("K32").Value = 407234.9
$ws.Range("M32").Value = -406947.9
$ws.Range("H61").Value = 3034.7334
$ws.Range("I61").Value = 2577.2222
$ws.Range("K61").Value = 2577.2222
$ws.Range("M61").Value = -2365.2222
$ws.Range("H102").Value = 2004.7778
$ws.Range("J102").Value = 2200
$ws.Range("L102").Value = 2200
$ws.Range("N102").Value = -5444
$ws.Range("H110").Value = 1246
$ws.Range("I110").Value = 1246
$ws.Range("K110").Value = 1246
$ws.Range("M110").Value = 799
$ws.Range("H116").Value = 2138.9375
$ws.Range("I116").Value = 2226.9167
$ws.Range("K116").Value = 2226.9167
$ws.Range("M116").Value = 67.08329999999978
$ws.Range("H132").Value = 5056.7095
$ws.Range("I132").Value = 4137.2856
$ws.Range("K132").Value = 12411.8568
$ws.Range("M132").Value = -9881.856800000001
$ws.Range("H136").Value = 3034.7334
$ws.Range("I136").Value = 2577.2222
$ws.Range("K136").Value = 7731.6666
$ws.Range("M136").Value = -5181.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2138.9375
$ws.Range("I3").Value = 2226.9167
$ws.Range("K3").Value = 2226.9167
$ws.Range("M3").Value = -2112.9167
$ws.Range("H86").Value = 90911680
$ws.Range("I86").Value = 90911680
$ws.Range("K86").Value = 90911680
$ws.Range("M86").Value = -90910557
$ws.Range("H89").Value = 90911680
$ws.Range("I89").Value = 90911680
$ws.Range("K89").Value = 454558400
$ws.Range("M89").Value = -454552784
$ws.Range("H94").Value = 1750
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1750
$ws.Range("N94").Value = -2652
$ws.Range("H99").Value = 885.7917
$ws.Range("I99").Value = 766.61536
$ws.Range("J99").Value = 1026.6364
$ws.Range("K99").Value = 766.61536
$ws.Range("L99").Value = 1026.6364
$ws.Range("M99").Value = 731.38464
$ws.Range("N99").Value = -4022.6364
$ws.Range("H105").Value = 10419856
$ws.Range("J105").Value = 3397.5
$ws.Range("L105").Value = 3397.5
$ws.Range("N105").Value = -6891.5
$ws.Range("H107").Value = 38101.074
$ws.Range("I107").Value = 50975.9
$ws.Range("J107").Value = 1315.8572
$ws.Range("K107").Value = 50975.9
$ws.Range("L107").Value = 1315.8572
$ws.Range("M107").Value = -49055.9
$ws.Range("N107").Value = -5155.8572
$ws.Range("M94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1489.5416
$ws.Range("I58").Value = 1232.8182
$ws.Range("J58").Value = 1706.7693
$ws.Range("K58").Value = 1232.8182
$ws.Range("L58").Value = 1706.7693
$ws.Range("M58").Value = -1029.8182
$ws.Range("N58").Value = -2112.7693
$ws.Range("H86").Value = 3223.2856
$ws.Range("I86").Value = 3427.9285
$ws.Range("J86").Value = 2814
$ws.Range("K86").Value = 3427.9285
$ws.Range("L86").Value = 2814
$ws.Range("M86").Value = -2304.9285
$ws.Range("N86").Value = -5060
$ws.Range("H89").Value = 3223.2856
$ws.Range("I89").Value = 3427.9285
$ws.Range("J89").Value = 2814
$ws.Range("K89").Value = 17139.6425
$ws.Range("L89").Value = 14070
$ws.Range("M89").Value = -11523.6425
$ws.Range("N89").Value = -25302
$ws.Range("H136").Value = 1489.5416
$ws.Range("I136").Value = 1232.8182
$ws.Range("J136").Value = 1706.7693
$ws.Range("K136").Value = 3698.4546
$ws.Range("L136").Value = 5120.3079
$ws.Range("M136").Value = -1148.4546
$ws.Range("N136").Value = -10220.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 990.84
$ws.Range("J113").Value = 1070.5238
$ws.Range("L113").Value = 3211.5714
$ws.Range("N113").Value = -7551.5714
$ws.Range("H131").Value = 1056.8572
$ws.Range("J131").Value = 1068.4584
$ws.Range("L131").Value = 3205.3752
$ws.Range("N131").Value = -13285.3752
$ws.Range("H132").Value = 2478.3635
$ws.Range("I132").Value = 2970
$ws.Range("J132").Value = 2333.7646
$ws.Range("K132").Value = 26730
$ws.Range("L132").Value = 21003.8814
$ws.Range("M132").Value = -24200
$ws.Range("N132").Value = -26063.8814
$ws.Range("H138").Value = 6292.222
$ws.Range("I138").Value = 882.25
$ws.Range("J138").Value = 7837.9287
$ws.Range("K138").Value = 2646.75
$ws.Range("L138").Value = 23513.7861
$ws.Range("M138").Value = 2493.25
$ws.Range("N138").Value = -33793.7861
$ws.Range("H140").Value = 1983.5333
$ws.Range("I140").Value = 1724.4445
$ws.Range("J140").Value = 2372.1667
$ws.Range("K140").Value = 5173.333500000001
$ws.Range("L140").Value = 7116.500100000001
$ws.Range("M140").Value = 6.66649999999936
$ws.Range("N140").Value = -17476.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 13867.375
$ws.Range("I93").Value = 15677
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 15677
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = -14429
$ws.Range("N93").Value = -3696
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H132").Value = 2675.9666
$ws.Range("I132").Value = 1657.7646
$ws.Range("J132").Value = 4007.4614
$ws.Range("K132").Value = 4973.293799999999
$ws.Range("L132").Value = 12022.3842
$ws.Range("M132").Value = -2443.293799999999
$ws.Range("N132").Value = -17082.3842
$ws.Range("H136").Value = 23812812
$ws.Range("I136").Value = 3001
$ws.Range("J136").Value = 55559224
$ws.Range("K136").Value = 9003
$ws.Range("L136").Value = 166677672
$ws.Range("M136").Value = -6453
$ws.Range("N136").Value = -166682772
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 77024.5
$ws.Range("I62").Value = 2049
$ws.Range("K62").Value = 2049
$ws.Range("M62").Value = -1425
$ws.Range("H65").Value = 77024.5
$ws.Range("I65").Value = 2049
$ws.Range("K65").Value = 10245
$ws.Range("M65").Value = -7125
$ws.Range("H96").Value = 4834.757
$ws.Range("I96").Value = 3185.75
$ws.Range("J96").Value = 5289.6553
$ws.Range("K96").Value = 3185.75
$ws.Range("L96").Value = 5289.6553
$ws.Range("M96").Value = -1812.75
$ws.Range("N96").Value = -8035.6553
$ws.Range("H132").Value = 4067249.5
$ws.Range("I132").Value = 2155.75
$ws.Range("J132").Value = 9806206
$ws.Range("K132").Value = 6467.25
$ws.Range("L132").Value = 29418618
$ws.Range("M132").Value = -3937.25
$ws.Range("N132").Value = -29423678
$ws.Range("H136").Value = 2384.795
$ws.Range("I136").Value = 2195.6296
$ws.Range("J136").Value = 2810.4167
$ws.Range("K136").Value = 6586.888800000001
$ws.Range("L136").Value = 8431.250100000001
$ws.Range("M136").Value = -4036.888800000001
$ws.Range("N136").Value = -13531.2501
